$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Helper cell used to stage apostrophe-prefixed (forced-text) values so that
# typing literal "false"/"true" strings doesn't get auto-converted to a
# native Excel boolean. We paste-special VALUES ONLY into the real target
# so the target's existing cell style/format is preserved, then remove the
# helper cell entirely.
$helper = $ws.Cells.Item(100, 26)

function Set-TextValue($cell, $text) {
    $helper.Value = "'" + $text
    $helper.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
    $excel.CutCopyMode = $false
    $helper.Delete()
}

# Row 7: Experimental = false
Set-TextValue $ws.Cells.Item(7, 2) "false"

# Row 8: Date value updated
$ws.Cells.Item(8, 2).Value = "2025-11-30T13:08:37+00:00"

# Row 17: Description text added
Set-TextValue $ws.Cells.Item(17, 2) "Codes for activity and recovery balance status"

Write-Host "B7 (Experimental): $($ws.Cells.Item(7,2).Value2)"
Write-Host "B8 (Date): $($ws.Cells.Item(8,2).Value2)"
Write-Host "B17 (Description): $($ws.Cells.Item(17,2).Value2)"
